{"js": "// Update the \"two-digit division\" answer table: each cell's\n// \"<dividend>\u00f7<divisor>=<quotient>, <remainder>\" text is replaced with a\n// newly generated problem/answer string. The mapping below reproduces the\n// diff exactly, in document order (old text -> new text).\nconst replacements = [\n  [\"64\u00f79=7, 1\", \"34\u00f79=3, 7\"],\n  [\"12\u00f78=1, 4\", \"38\u00f77=5, 3\"],\n  [\"82\u00f74=20, 2\", \"80\u00f79=8, 8\"],\n  [\"57\u00f79=6, 3\", \"22\u00f76=3, 4\"],\n  [\"79\u00f77=11, 2\", \"63\u00f76=10, 3\"],\n  [\"93\u00f72=46, 1\", \"65\u00f76=10, 5\"],\n  [\"52\u00f73=17, 1\", \"87\u00f76=14, 3\"],\n  [\"38\u00f76=6, 2\", \"93\u00f75=18, 3\"],\n  [\"61\u00f73=20, 1\", \"11\u00f76=1, 5\"],\n  [\"55\u00f76=9, 1\", \"63\u00f78=7, 7\"],\n  [\"81\u00f78=10, 1\", \"30\u00f79=3, 3\"],\n  [\"70\u00f74=17, 2\", \"71\u00f74=17, 3\"],\n  [\"52\u00f78=6, 4\", \"50\u00f79=5, 5\"],\n  [\"66\u00f77=9, 3\", \"86\u00f76=14, 2\"],\n  [\"75\u00f72=37, 1\", \"47\u00f72=23, 1\"],\n  [\"41\u00f75=8, 1\", \"74\u00f76=12, 2\"],\n  [\"42\u00f75=8, 2\", \"26\u00f74=6, 2\"],\n  [\"93\u00f77=13, 2\", \"29\u00f78=3, 5\"],\n  [\"73\u00f77=10, 3\", \"91\u00f73=30, 1\"],\n  [\"17\u00f78=2, 1\", \"74\u00f76=12, 2\"],\n  [\"62\u00f74=15, 2\", \"93\u00f76=15, 3\"],\n  [\"85\u00f75=17, 0\", \"30\u00f79=3, 3\"],\n  [\"50\u00f74=12, 2\", \"79\u00f77=11, 2\"],\n  [\"34\u00f72=17, 0\", \"30\u00f75=6, 0\"],\n  [\"43\u00f73=14, 1\", \"71\u00f75=14, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  // Each source string is unique in the document, so replace every match\n  // found (expected to be exactly one) with the new text.\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the \"two-digit division\" answer table: each cell's\n# \"<dividend>\u00f7<divisor>=<quotient>, <remainder>\" text is replaced with a\n# newly generated problem/answer string. The mapping below reproduces the\n# diff exactly, in document order (old text -> new text).\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"64\u00f79=7, 1\", \"34\u00f79=3, 7\"),\n    @(\"12\u00f78=1, 4\", \"38\u00f77=5, 3\"),\n    @(\"82\u00f74=20, 2\", \"80\u00f79=8, 8\"),\n    @(\"57\u00f79=6, 3\", \"22\u00f76=3, 4\"),\n    @(\"79\u00f77=11, 2\", \"63\u00f76=10, 3\"),\n    @(\"93\u00f72=46, 1\", \"65\u00f76=10, 5\"),\n    @(\"52\u00f73=17, 1\", \"87\u00f76=14, 3\"),\n    @(\"38\u00f76=6, 2\", \"93\u00f75=18, 3\"),\n    @(\"61\u00f73=20, 1\", \"11\u00f76=1, 5\"),\n    @(\"55\u00f76=9, 1\", \"63\u00f78=7, 7\"),\n    @(\"81\u00f78=10, 1\", \"30\u00f79=3, 3\"),\n    @(\"70\u00f74=17, 2\", \"71\u00f74=17, 3\"),\n    @(\"52\u00f78=6, 4\", \"50\u00f79=5, 5\"),\n    @(\"66\u00f77=9, 3\", \"86\u00f76=14, 2\"),\n    @(\"75\u00f72=37, 1\", \"47\u00f72=23, 1\"),\n    @(\"41\u00f75=8, 1\", \"74\u00f76=12, 2\"),\n    @(\"42\u00f75=8, 2\", \"26\u00f74=6, 2\"),\n    @(\"93\u00f77=13, 2\", \"29\u00f78=3, 5\"),\n    @(\"73\u00f77=10, 3\", \"91\u00f73=30, 1\"),\n    @(\"17\u00f78=2, 1\", \"74\u00f76=12, 2\"),\n    @(\"62\u00f74=15, 2\", \"93\u00f76=15, 3\"),\n    @(\"85\u00f75=17, 0\", \"30\u00f79=3, 3\"),\n    @(\"50\u00f74=12, 2\", \"79\u00f77=11, 2\"),\n    @(\"34\u00f72=17, 0\", \"30\u00f75=6, 0\"),\n    @(\"43\u00f73=14, 1\", \"71\u00f75=14, 1\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $found = $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $old\"\n    }\n}\n"}
